$d = $word.ActiveDocument

# Mapping of old text -> new text (date line + 24 multiplication problems).
# Order mirrors document order; each old value is unique in the document so
# a simple sequential Find/Replace (no wildcards) is unambiguous and safe.
$pairs = @(
    @("2025-06-11 Wednesday", "2025-06-12 Thursday"),
    @("588×2=", "364×9="),
    @("965×4=", "569×3="),
    @("653×6=", "764×7="),
    @("414×5=", "151×4="),
    @("353×8=", "954×3="),
    @("472×9=", "933×8="),
    @("668×6=", "846×4="),
    @("824×8=", "400×5="),
    @("558×2=", "379×6="),
    @("539×6=", "693×8="),
    @("917×6=", "843×9="),
    @("917×9=", "257×6="),
    @("860×8=", "942×7="),
    @("602×8=", "569×3="),
    @("880×8=", "759×3="),
    @("241×3=", "642×2="),
    @("745×8=", "499×6="),
    @("101×5=", "578×3="),
    @("720×4=", "249×4="),
    @("755×6=", "937×8="),
    @("691×6=", "554×8="),
    @("509×8=", "136×2="),
    @("465×7=", "882×9="),
    @("342×2=", "529×6="),
    @("860×9=", "724×7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
